# Insert a new weekly price record at row 218 (Betarraga, Feria Lagunitas de
# Puerto Montt), pushing the existing rows 218..285 down to 219..286.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(218).Insert()

$ws.Range("A218").Value = 4
$ws.Range("B218").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C218").Value = "Los Lagos"
$ws.Range("D218").Value = 44663
$ws.Range("E218").Value = 10
$ws.Range("F218").Value = 100114014
$ws.Range("G218").Value = "Betarraga"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 800
$ws.Range("K218").Value = 1000
$ws.Range("L218").Value = 1000
$ws.Range("M218").Value = 1000
$ws.Range("N218").Value = "$/paquete 5 unidades"
$ws.Range("O218").Value = "Región del Maule"
$ws.Range("P218").Value = 200
$ws.Range("Q218").Value = 5
$ws.Range("R218").Value = "Hortaliza"
